# Weekly update: insert a new week's worth of data (2 rows: "Primera" and
# "Segunda" quality) at row 508, pushing all subsequent rows down by 2.
# This mirrors the commit "Fruta / hortaliza, semanal" (weekly fruit/veg
# price refresh) where the newest observation is inserted in date order
# and older rows shift down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 508:509, shifting existing rows (old 508..608)
# down to (510..610).
$ws.Rows("508:509").Insert()

# New row 508 - Calidad "Primera"
$ws.Range("A508").Value = 11
$ws.Range("B508").Value = "Vega Monumental Concepción"
$ws.Range("C508").Value = "Bíobío"
$ws.Range("D508").Value = 45258
$ws.Range("E508").Value = 8
$ws.Range("F508").Value = 100114014
$ws.Range("G508").Value = "Betarraga"
$ws.Range("H508").Value = "Sin especificar"
$ws.Range("I508").Value = "Primera"
$ws.Range("J508").Value = 650
$ws.Range("K508").Value = 600
$ws.Range("L508").Value = 650
$ws.Range("M508").Value = 627
$ws.Range("N508").Value = "$/paquete 5 unidades"
$ws.Range("O508").Value = "Región Metropolitana"
$ws.Range("P508").Value = 125
$ws.Range("Q508").Value = 5
$ws.Range("R508").Value = "Hortaliza"

# New row 509 - Calidad "Segunda"
$ws.Range("A509").Value = 11
$ws.Range("B509").Value = "Vega Monumental Concepción"
$ws.Range("C509").Value = "Bíobío"
$ws.Range("D509").Value = 45258
$ws.Range("E509").Value = 8
$ws.Range("F509").Value = 100114014
$ws.Range("G509").Value = "Betarraga"
$ws.Range("H509").Value = "Sin especificar"
$ws.Range("I509").Value = "Segunda"
$ws.Range("J509").Value = 200
$ws.Range("K509").Value = 500
$ws.Range("L509").Value = 500
$ws.Range("M509").Value = 500
$ws.Range("N509").Value = "$/paquete 5 unidades"
$ws.Range("O509").Value = "Región Metropolitana"
$ws.Range("P509").Value = 100
$ws.Range("Q509").Value = 5
$ws.Range("R509").Value = "Hortaliza"
